# Add a "Category" header in A1 (matching the existing header style used by B1:W1),
# and reset the style of A2:A46 (the category label column) back to the default/Normal
# style, since that formatting moves up into the new header row instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the new header cell A1 with the same formatting as the other header cells.
$ws.Range("A1").Value = "Category"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# 2) Remove the (now redundant) header-style formatting from A2:A46, restoring it
#    to the workbook's default "Normal" style.
$ws.Range("A2:A46").Style = "Normal"
